# DiscountRules.xlsx: update row 24
#  - A24 becomes a new text value "TEST LAST" (replacing the previous blank entry)
#  - E24 is cleared back to blank (it previously held "5x5 0048")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A24").Value = "TEST LAST"

$ws.Range("E24").ClearContents()
$ws.Range("E24").Style = "Normal"
